$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 previously held "Pastor Rev Heidi" / "pastorheidi@wfcrc.ca" / "345-9893".
# Update it to "Ken M" / "Ken_Mullins@sil.org" (as a mailto hyperlink) / "345-9893" (unchanged).
$ws.Range("A11").Value = "Ken M"

$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:Ken_Mullins@sil.org", "", "", "Ken_Mullins@sil.org") | Out-Null

$ws.Range("B11").Select()
